# Applies the "quiz marksheet" regrading fix:
#  - Row 10-12 summary counters & labels updated (Right/Wrong/NotAttempt/Max, Marking, Total)
#  - C11 marking-for-wrong-answer is stored as a real number instead of text ("-1" -> -1)
#  - Student answers are now filled into the "Student Ans" column(s) so the
#    correct / incorrect / not-attempted styling is meaningful, instead of the
#    column being left blank for (almost) every question
#  - The third answer block (columns G/H) is removed entirely, and the second
#    block (columns D/E) is trimmed down to only the first few rows
#  - Used range shrinks from A5:H40 down to A5:E40 as a consequence

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary block (rows 10-12)
# ---------------------------------------------------------------------------

# Labels in column A now use the same "title" style as the row-9 header
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

$ws.Range("A10").Value = "No."
$ws.Range("A11").Value = "Marking"
$ws.Range("A12").Value = "Total"

# Right / Wrong / Not Attempt / Max counts
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

# Marks awarded per right / wrong answer (must be numeric, not text, so math
# downstream keeps working even with negative / float marking schemes)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Totals: right marks, wrong marks, final score/out-of-max
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "57/112"

# ---------------------------------------------------------------------------
# Third answer block (columns G/H) is dropped completely
# ---------------------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------------
# Second answer block (columns D/E): only rows 16-18 survive
# ---------------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

# D17 stays blank/not-attempted (style unchanged), only E17 keeps its value

$ws.Range("D18").Style = "incorrectStyle"
$ws.Range("D18").Value = "Option B"

# ---------------------------------------------------------------------------
# First answer block (columns A/B), rows 16-40: fill in the student's answer
# ---------------------------------------------------------------------------

function Set-Answer($cellRef, $styleName, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Style = $styleName
    $rng.Value = $value
}

Set-Answer "A16" "correctStyle"   "Option A"
Set-Answer "A17" "correctStyle"   "Option D"
Set-Answer "A18" "correctStyle"   "Option B"
Set-Answer "A19" "correctStyle"   "Option C"
Set-Answer "A20" "correctStyle"   "Option B"
Set-Answer "A21" "correctStyle"   "Option C"
# A22, A23 stay blank/not-attempted
Set-Answer "A24" "correctStyle"   "Option A"
Set-Answer "A25" "correctStyle"   "Option A"
Set-Answer "A26" "incorrectStyle" "Option B"
Set-Answer "A27" "correctStyle"   "Option A"
Set-Answer "A28" "correctStyle"   "Option D"
# A29 stays blank/not-attempted
Set-Answer "A30" "incorrectStyle" "Option C"
# A31 stays blank/not-attempted
Set-Answer "A32" "correctStyle"   "Option C"
Set-Answer "A33" "correctStyle"   "Option D"
# A34, A35 stay blank/not-attempted
Set-Answer "A36" "correctStyle"   "Option A"
# A37, A38 stay blank/not-attempted
Set-Answer "A39" "correctStyle"   "Option D"
# A40 stays blank/not-attempted
